$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D11:F11").NumberFormat = "@"

$ws.Range("A11").Value = "2P"
$ws.Range("B11").Value = "Papelería"
$ws.Range("C11").Value = "Resma x 100 hojas"
$ws.Range("D11").Value = "10"
$ws.Range("E11").Value = "10000"
$ws.Range("F11").Value = "15000"
$ws.Range("G11").Value = 45818.96507344174

$ws.Range("D11:F11").Style = "Normal"
